# Update column F (dSF) values with repulled data, per the commit:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    4  = -1
    6  = -2
    8  = -2
    9  = -4
    10 = 4
    11 = -4
    12 = -2
    13 = 0
    15 = 4
    16 = -4
    17 = -1
    18 = 2
    19 = 4
    20 = 2
    21 = 1
    22 = 3
    23 = -3
    24 = 3
    25 = 3
    26 = 1
    27 = 3
    28 = 2
    29 = 1
    30 = 1
    32 = -1
    34 = 7
    35 = -1
    36 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
